# Applies three small word-level edits inside the "Project Report" document:
#   1. "This project involves all 3 parts ..."      -> "... all three parts ..."
#   2. "We loaded the 2 CSV files into ..."          -> "... the two CSV files ..."
#   3. "... using progress query ..."                -> "... using SQL query ..."
#
# Each edit is performed by first locating the short, unambiguous phrase that
# contains the word to change, then narrowing the search to just that one
# word and overtyping it (matching how a person would select the word in
# Word and type the replacement) so the surrounding run formatting
# (Arial, bold, 28 half-points) is preserved untouched.

$d = $word.ActiveDocument

function Replace-WordInPhrase($phrase, $oldWord, $newWord) {
    $outer = $d.Content
    $foundOuter = $outer.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $foundOuter) {
        throw "Could not locate phrase: $phrase"
    }

    $inner = $d.Range($outer.Start, $outer.End)
    $foundInner = $inner.Find.Execute($oldWord, $true, $false, $false, $false, $false, $true, 1, $false, $newWord, 2)
    if (-not $foundInner) {
        throw "Could not locate word '$oldWord' inside phrase: $phrase"
    }
}

# 1) "all 3 parts" -> "all three parts"
Replace-WordInPhrase "all 3 parts of the ETL process" "3" "three"

# 2) "loaded the 2 CSV files" -> "loaded the two CSV files"
Replace-WordInPhrase "loaded the 2 CSV files" "2" "two"

# 3) "using progress query" -> "using SQL query"
Replace-WordInPhrase "using progress query" "progress" "SQL"
